# This script reproduces the commit's effect: the first data row (the
# "2024-06-16 合肥·人气COSER次元会（免费展）" event, previously row 2) was
# removed from the "展览" (Exhibition) sheet and the "全部类型" (All types)
# sheet. All subsequent rows shifted up by one, and the running index in
# column A was renumbered sequentially again (1, 2, 3, ...).

$wb = $excel.ActiveWorkbook

# Sheet 1 = 展览 (Exhibition) -- originally A1:I22, becomes A1:I21
$wsExhibition = $wb.Worksheets.Item(1)
$wsExhibition.Rows.Item(2).Delete()
$lastRow = $wsExhibition.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $wsExhibition.Cells.Item($r, 1).Value2 = $r - 1
}

# Sheet 4 = 全部类型 (All types) -- originally A1:I23, becomes A1:I22
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Rows.Item(2).Delete()
$lastRow = $wsAll.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $wsAll.Cells.Item($r, 1).Value2 = $r - 1
}
